# Updated cryptos list (price + 1h volume change) as per the latest scrape.
# Columns: D = Price (text, dot-grouped), E = Volume(1h) change (text, "  +x.xx%  ").
# D/E were stored as plain text in the sheet, so any value that *looks* like a
# bare number (e.g. "0.9995", "238.04") must be forced back to Text first -
# otherwise Excel's COM layer auto-coerces it to a numeric cell on assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; NewD = "25.535.52"; NewE = "  +2.56%  "; ForceText = $False },
    @{ Row = 3; NewD = "1.669.94"; NewE = "  +2.12%  "; ForceText = $False },
    @{ Row = 4; NewD = "0.9995"; NewE = "  +0.22%  "; ForceText = $True },
    @{ Row = 5; NewD = "238.04"; NewE = "  +0.96%  "; ForceText = $True },
    @{ Row = 6; NewD = "1.001"; NewE = "  -0.04%  "; ForceText = $True },
    @{ Row = 7; NewD = $null; NewE = "  +1.40%  "; ForceText = $False },
    @{ Row = 8; NewD = "0.2617"; NewE = "  +3.03%  "; ForceText = $True },
    @{ Row = 9; NewD = "0.06172"; NewE = "  +2.96%  "; ForceText = $True },
    @{ Row = 10; NewD = "1.671.83"; NewE = "  +2.09%  "; ForceText = $False },
    @{ Row = 11; NewD = "0.06973"; NewE = $null; ForceText = $True },
    @{ Row = 12; NewD = "14.81"; NewE = "  +0.59%  "; ForceText = $True },
    @{ Row = 13; NewD = "0.5875"; NewE = "  -4.22%  "; ForceText = $True },
    @{ Row = 14; NewD = "4.376"; NewE = $null; ForceText = $True },
    @{ Row = 15; NewD = "75.22"; NewE = "  +3.69%  "; ForceText = $True },
    @{ Row = 16; NewD = $null; NewE = "  -0.02%  "; ForceText = $False },
    @{ Row = 17; NewD = "1.000"; NewE = "  +0.12%  "; ForceText = $True },
    @{ Row = 18; NewD = "25.533.82"; NewE = "  +2.51%  "; ForceText = $False },
    @{ Row = 19; NewD = "0.000006747"; NewE = "  +2.91%  "; ForceText = $True },
    @{ Row = 20; NewD = "11.44"; NewE = "  +3.34%  "; ForceText = $True },
    @{ Row = 21; NewD = "1.887.62"; NewE = "  +2.29%  "; ForceText = $False },
    @{ Row = 22; NewD = "4.450"; NewE = "  +2.14%  "; ForceText = $True },
    @{ Row = 23; NewD = "8.799"; NewE = "  +2.75%  "; ForceText = $True },
    @{ Row = 24; NewD = "5.266"; NewE = "  +0.34%  "; ForceText = $True },
    @{ Row = 25; NewD = "136.47"; NewE = "  +2.25%  "; ForceText = $True },
    @{ Row = 26; NewD = "15.03"; NewE = "  +1.70%  "; ForceText = $True },
    @{ Row = 27; NewD = $null; NewE = "  +1.54%  "; ForceText = $False },
    @{ Row = 28; NewD = "1.724"; NewE = "  +4.47%  "; ForceText = $True },
    @{ Row = 29; NewD = "104.59"; NewE = "  +2.16%  "; ForceText = $True },
    @{ Row = 30; NewD = "3.977"; NewE = "  +6.29%  "; ForceText = $True },
    @{ Row = 31; NewD = "0.07870"; NewE = "  +2.07%  "; ForceText = $True },
    @{ Row = 32; NewD = "3.623"; NewE = "  +2.04%  "; ForceText = $True },
    @{ Row = 33; NewD = "0.9995"; NewE = $null; ForceText = $True },
    @{ Row = 34; NewD = "0.04270"; NewE = "  -0.56%  "; ForceText = $True },
    @{ Row = 35; NewD = "2.619"; NewE = "  +0.85%  "; ForceText = $True },
    @{ Row = 36; NewD = "0.9538"; NewE = "  +4.15%  "; ForceText = $True },
    @{ Row = 37; NewD = "0.6070"; NewE = "  +5.17%  "; ForceText = $True },
    @{ Row = 38; NewD = $null; NewE = "  +1.29%  "; ForceText = $False },
    @{ Row = 39; NewD = "0.8979"; NewE = "  +9.23%  "; ForceText = $True },
    @{ Row = 40; NewD = $null; NewE = "  +0.17%  "; ForceText = $False },
    @{ Row = 41; NewD = "1.855"; NewE = "  +3.61%  "; ForceText = $True },
    @{ Row = 42; NewD = "0.01478"; NewE = "  -4.31%  "; ForceText = $True },
    @{ Row = 43; NewD = "96.18"; NewE = "  -0.79%  "; ForceText = $True },
    @{ Row = 44; NewD = "0.3754"; NewE = "  +1.61%  "; ForceText = $True },
    @{ Row = 45; NewD = "4.888"; NewE = "  +3.45%  "; ForceText = $True },
    @{ Row = 46; NewD = "0.1118"; NewE = "  +1.95%  "; ForceText = $True },
    @{ Row = 47; NewD = "6.210"; NewE = "  +2.66%  "; ForceText = $True },
    @{ Row = 48; NewD = "0.05264"; NewE = "  +1.06%  "; ForceText = $True },
    @{ Row = 49; NewD = "29.96"; NewE = "  +1.77%  "; ForceText = $True },
    @{ Row = 50; NewD = "7.438"; NewE = "  +4.34%  "; ForceText = $True },
    @{ Row = 51; NewD = $null; NewE = "  +0.15%  "; ForceText = $False }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.NewD) {
        $dCell = $ws.Cells.Item($r, 4)
        if ($u.ForceText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.NewD
    }
    if ($null -ne $u.NewE) {
        $ws.Cells.Item($r, 5).Value = $u.NewE
    }
}
